$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "current_task"
$ws.Range("B2").Value = "str"
$ws.Range("A3").Value = "time:timestamp"
$ws.Range("B3").Value = "datetime"
$ws.Range("A4").Value = "complete_service_time"
$ws.Range("B4").Value = "str"
$ws.Range("A5").Value = "SubProcessID"
$ws.Range("B5").Value = "str"
$ws.Range("A6").Value = "response_status_code"
$ws.Range("B6").Value = "float"
$ws.Range("A7").Value = "human_workstation_green_button_pressed"
$ws.Range("B7").Value = "float"
$ws.Range("A8").Value = "org:resource"
$ws.Range("B8").Value = "str"
$ws.Range("A9").Value = "identifier:id"
$ws.Range("B9").Value = "str"
$ws.Range("A10").Value = "process_model_id"
$ws.Range("B10").Value = "str"
$ws.Range("A11").Value = "lifecycle:state"
$ws.Range("B11").Value = "str"
$ws.Range("A12").Value = "lifecycle:transition"
$ws.Range("B12").Value = "str"
$ws.Range("A13").Value = "case"
$ws.Range("B13").Value = "str"
$ws.Range("A14").Value = "parameters"
$ws.Range("B14").Value = "dict"
$ws.Range("A15").Value = "operation_end_time"
$ws.Range("B15").Value = "datetime"
$ws.Range("A16").Value = "unsatisfied_condition_description"
$ws.Range("B16").Value = "str"
$ws.Range("A17").Value = "event_id"
$ws.Range("B17").Value = "str"
$ws.Range("A18").Value = "case:concept:name"
$ws.Range("B18").Value = "str"
$ws.Range("A19").Value = "planned_operation_time"
$ws.Range("B19").Value = "str"
$ws.Range("A20").Value = "requested_service_url"
$ws.Range("B20").Value = "str"
$ws.Range("A21").Value = "concept:name"
$ws.Range("B21").Value = "str"
